$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 46063
$ws.Range("B2").Value = 0.57
$ws.Range("C2").Value = 0.05
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.01
$ws.Range("I2").Value = 0.49
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2.46
$ws.Range("L2").Value = 2.44
$ws.Range("M2").Value = 0.8100000000000001
$ws.Range("N2").Value = 0.31
$ws.Range("O2").Value = 0.07000000000000001
$ws.Range("P2").Value = 0.02
$ws.Range("Q2").Value = 0.07000000000000001
$ws.Range("R2").Value = 0.5
$ws.Range("S2").Value = 2.03
$ws.Range("T2").Value = 2.04
$ws.Range("U2").Value = 5.6
$ws.Range("V2").Value = 11.55
$ws.Range("W2").Value = 10.06
$ws.Range("X2").Value = 2.6
$ws.Range("Y2").Value = 0.95
$ws.Range("Z2").Value = 1.82
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 6.29
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 10.8
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 3.82
$ws.Range("AG2").Value = "0h-23h"
